# "include references to indicators"
# The indicator codes stored as shared strings were reformatted from
# dotted notation (e.g. "1.1.a") to space-separated notation (e.g. "1.1 a").
# These values live in column C of Sheet1, rows 11-13 and 29-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "1.1 a"
$ws.Range("C12").Value = "1.4 a"
$ws.Range("C13").Value = "1.6 c"
$ws.Range("C29").Value = "1.1 c"
$ws.Range("C30").Value = "1.2 b"

# Match the author's final cursor position recorded in the sheet view.
$ws.Range("C31").Select()
